$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new expense rows (Pos EAD related purchases)
$ws.Range("A14").Value = "Placa de Video "
$ws.Range("B14").Value = 80

$ws.Range("A15").Value = "Adaptador Dvi / Vga"
$ws.Range("B15").Value = 20.8

# Move/update the active selection to where the user left off entering data
$ws.Range("B16").Select()

$wb.Save()
